# Tidsrapportering.xlsx - "Updated the time report"
#
# Marcus logged additional hours for week 18 (column G) and week 19
# (column H) on the "Kodning" row, which ripples through the weekly
# total (K4), the per-week total row (row 12) and the grand total (C15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marcus")

# New hours entered for "Kodning": 11h in week 18 (col G), 2h in week 19 (col H)
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 2

# The week-18 column total (row 12) is re-entered explicitly so it only
# sums the data rows (4:11), breaking it out of the shared formula that
# used to start at row 3.
$ws.Range("G12").Formula = "=SUM(G4:G11)"

# Leave the selection on the last cell that was edited.
$ws.Range("H4").Select()
